$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.558.90"
$ws.Range("E2").Value = "  -2.56%  "

$ws.Range("D3").Value = "2.412.38"
$ws.Range("E3").Value = "  -2.20%  "

$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'564.79"
$ws.Range("E5").Value = "  -3.26%  "

$ws.Range("D6").Value = "'137.19"
$ws.Range("E6").Value = "  -3.89%  "

$ws.Range("E7").Value = "  +0.32%  "

$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("D9").Value = "2.396.45"
$ws.Range("E9").Value = "  -2.72%  "

$ws.Range("E10").Value = "  -5.78%  "

$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").Value = "'5.04"
$ws.Range("E12").Value = "  -3.11%  "

$ws.Range("E13").Value = "  -1.98%  "

$ws.Range("D14").Value = "'25.64"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "2.811.12"
$ws.Range("E15").Value = "  -2.99%  "

$ws.Range("D16").Value = "'0.0000166"
$ws.Range("E16").Value = "  -4.16%  "

$ws.Range("D17").Value = "60.836.42"
$ws.Range("E17").Value = "  -1.91%  "

$ws.Range("D18").Value = "2.390.75"
$ws.Range("E18").Value = "  -2.80%  "

$ws.Range("D19").Value = "'8.08"
$ws.Range("E19").Value = "  +9.98%  "

$ws.Range("D20").Value = "'10.51"
$ws.Range("E20").Value = "  -1.95%  "

$ws.Range("D21").Value = "'322.32"
$ws.Range("E21").Value = "  -1.50%  "

$ws.Range("E22").Value = "  -1.59%  "

$ws.Range("D23").Value = "'6.15"
$ws.Range("E23").Value = "  -1.37%  "

$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").Value = "'1.79"
$ws.Range("E25").Value = "  -8.14%  "

$ws.Range("D26").Value = "'64.06"
$ws.Range("E26").Value = "  -2.05%  "

$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").Value = "'556.58"
$ws.Range("E27").Value = "  -5.34%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'8.09"
$ws.Range("E28").Value = "  -12.82%  "

$ws.Range("D29").Value = "2.523.67"
$ws.Range("E29").Value = "  -1.78%  "

$ws.Range("D30").Value = "0.0₃0909"
$ws.Range("E30").Value = "  -3.59%  "

$ws.Range("D31").Value = "'7.86"
$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").Value = "'1.29"
$ws.Range("E32").Value = "  -7.61%  "

$ws.Range("E33").Value = "  -5.24%  "

$ws.Range("E34").Value = "  -2.76%  "

$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'153.64"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.41"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("D38").Value = "'0.367"
$ws.Range("E38").Value = "  -1.83%  "

$ws.Range("D39").Value = "'4.50"
$ws.Range("E39").Value = "  -6.20%  "

$ws.Range("D40").Value = "'18.14"
$ws.Range("E40").Value = "  -1.43%  "

$ws.Range("D41").Value = "'5.05"
$ws.Range("E41").Value = "  -3.13%  "

$ws.Range("E42").Value = "  +0.03%  "

$ws.Range("E43").Value = "  -3.82%  "

$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").Value = "0.0₆0289"
$ws.Range("E44").Value = "  -1.92%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.28"
$ws.Range("E45").Value = "  -5.26%  "

$ws.Range("D46").Value = "'142.28"
$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("D47").Value = "'3.49"
$ws.Range("E47").Value = "  -2.71%  "

$ws.Range("D48").Value = "'0.582"
$ws.Range("E48").Value = "  -3.59%  "

$ws.Range("D49").Value = "'0.0496"
$ws.Range("E49").Value = "  -3.67%  "

$ws.Range("D50").Value = "'18.97"
$ws.Range("E50").Value = "  -5.07%  "

$ws.Range("D51").Value = "'0.0893"
$ws.Range("E51").Value = "  -0.74%  "
